$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C to fit the new, longer activity descriptions.
$ws.Columns.Item(3).ColumnWidth = 41.14

# New "Dauer (berechnet)" helper cell in E2, formatted as [h]:mm.
$ws.Range("E2").NumberFormat = "[h]:mm"

$ws.Range("A48").Value = 44310
$ws.Range("B48").Value = 0.10416666666666667
$ws.Range("C48").Value = "Implementierung"
$ws.Range("D48").Value = "Tests"
$ws.Range("A49").Value = 44311
$ws.Range("B49").Value = 0.0625
$ws.Range("C49").Value = "Koordination und Projektmanagement"
$ws.Range("D49").Value = "Teilnahme an Meeting"
$ws.Range("A50").Value = 44311
$ws.Range("B50").Value = 0.10416666666666667
$ws.Range("C50").Value = "Implementierung"
$ws.Range("D50").Value = "User Case Game Erstellen"
$ws.Range("A51").Value = 44312
$ws.Range("B51").Value = 0.08333333333333333
$ws.Range("C51").Value = "LV-Einheit"
$ws.Range("D51").Value = "Teilnahme an PS-Einheit"
$ws.Range("A52").Value = 44313
$ws.Range("B52").Value = 0.10416666666666667
$ws.Range("C52").Value = "Implementierung"
$ws.Range("D52").Value = "User Case Teams Erstellen"
$ws.Range("A53").Value = 44313
$ws.Range("B53").Value = 0.0625
$ws.Range("C53").Value = "Koordination und Projektmanagement"
$ws.Range("D53").Value = "Teilnahme an Meeting"
$ws.Range("A54").Value = 44315
$ws.Range("B54").Value = 0.125
$ws.Range("C54").Value = "Implementierung"
$ws.Range("D54").Value = "Game- and Team-related Tests"
$ws.Range("A55").Value = 44316
$ws.Range("B55").Value = 0.10416666666666667
$ws.Range("C55").Value = "Implementierung"
$ws.Range("D55").Value = "Use Case Term/Topic-Upload"
$ws.Range("A56").Value = 44317
$ws.Range("B56").Value = 0.08333333333333333
$ws.Range("C56").Value = "Implementierung"
$ws.Range("D56").Value = "Term- and Topic-related Tests"
$ws.Range("A57").Value = 44317
$ws.Range("B57").Value = 0.020833333333333332
$ws.Range("C57").Value = "Koordination und Projektmanagement"
$ws.Range("D57").Value = "Teilnahme an Meeting"
$ws.Range("A58").Value = 44318
$ws.Range("B58").Value = 0.08333333333333333
$ws.Range("C58").Value = "Implementierung"
$ws.Range("D58").Value = "Feature Einheitliche VirtualLobby"
$ws.Range("A59").Value = 44318
$ws.Range("B59").Value = 0.041666666666666664
$ws.Range("C59").Value = "Koordination und Projektmanagement"
$ws.Range("D59").Value = "Teilnahme an Meeting"
$ws.Range("A60").Value = 44320
$ws.Range("B60").Value = 0.0625
$ws.Range("C60").Value = "Implementierung"
$ws.Range("D60").Value = "Use Case Spielereinladung"
$ws.Range("A61").Value = 44320
$ws.Range("B61").Value = 0.041666666666666664
$ws.Range("C61").Value = "Koordination und Projektmanagement"
$ws.Range("D61").Value = "Teilnahme an Meeting"
$ws.Range("A62").Value = 44322
$ws.Range("B62").Value = 0.0625
$ws.Range("C62").Value = "Implementierung"
$ws.Range("D62").Value = "Tests"
$ws.Range("A63").Value = 44323
$ws.Range("B63").Value = 0.125
$ws.Range("C63").Value = "Implementierung"
$ws.Range("D63").Value = "Use Case Spielereinladung"
$ws.Range("A64").Value = 44324
$ws.Range("B64").Value = 0.08333333333333333
$ws.Range("C64").Value = "Implementierung"
$ws.Range("D64").Value = "Erweiterung/Verbesserung bestehender Funktionalitäten"
$ws.Range("A65").Value = 44324
$ws.Range("B65").Value = 0.041666666666666664
$ws.Range("C65").Value = "Koordination und Projektmanagement"
$ws.Range("D65").Value = "Teilnahme an Meeting"
$ws.Range("A66").Value = 44325
$ws.Range("B66").Value = 0.08333333333333333
$ws.Range("C66").Value = "Implementierung"
$ws.Range("D66").Value = "Feature Anwesenheitscheck im Gameroom"
$ws.Range("A67").Value = 44325
$ws.Range("B67").Value = 0.041666666666666664
$ws.Range("C67").Value = "Koordination und Projektmanagement"
$ws.Range("D67").Value = "Teilnahme an Meeting"
$ws.Range("A68").Value = 44326
$ws.Range("B68").Value = 0.08333333333333333
$ws.Range("C68").Value = "LV-Einheit"
$ws.Range("D68").Value = "Teilnahme an PS-Einheit"
$ws.Range("A69").Value = 44326
$ws.Range("B69").Value = 0.08333333333333333
$ws.Range("C69").Value = "Implementierung"
$ws.Range("D69").Value = "Feature Anwesenheitscheck im Gameroom"
$ws.Range("A70").Value = 44327
$ws.Range("B70").Value = 0.125
$ws.Range("C70").Value = "Implementierung"
$ws.Range("D70").Value = "Feature Echtzeitupdate"
$ws.Range("A71").Value = 44327
$ws.Range("B71").Value = 0.041666666666666664
$ws.Range("C71").Value = "Koordination und Projektmanagement"
$ws.Range("D71").Value = "Teilnahme an Meeting"
$ws.Range("A72").Value = 44328
$ws.Range("B72").Value = 0.125
$ws.Range("C72").Value = "Implementierung"
$ws.Range("D72").Value = "Feature Echtzeitupdate"
$ws.Range("A73").Value = 44329
$ws.Range("B73").Value = 0.125
$ws.Range("C73").Value = "Implementierung"
$ws.Range("D73").Value = "Datenbanktrennung und Fixes"
$ws.Range("A74").Value = 44329
$ws.Range("B74").Value = 0.25
$ws.Range("C74").Value = "Koordination und Projektmanagement"
$ws.Range("D74").Value = "Teilnahme an Meeting"
$ws.Range("A75").Value = 44333
$ws.Range("B75").Value = 0.041666666666666664
$ws.Range("C75").Value = "Systemtest (fremdes System)"
$ws.Range("D75").Value = "Teilnahme an Meeting"
$ws.Range("A76").Value = 44334
$ws.Range("B76").Value = 0.0625
$ws.Range("C76").Value = "Systemtest (fremdes System)"
$ws.Range("D76").Value = "User-related Tests"
$ws.Range("A77").Value = 44335
$ws.Range("B77").Value = 0.08333333333333333
$ws.Range("C77").Value = "Systemtest (fremdes System)"
$ws.Range("D77").Value = "Teilnahme an Meeting"
$ws.Range("A78").Value = 44336
$ws.Range("B78").Value = 0.125
$ws.Range("C78").Value = "Systemtest (fremdes System)"
$ws.Range("D78").Value = "Teilnahme an Meeting"
$ws.Range("A79").Value = 44343
$ws.Range("B79").Value = 0.125
$ws.Range("C79").Value = "Implementierung"
$ws.Range("D79").Value = "Bugfixes Abnahmetest"
$ws.Range("A80").Value = 44343
$ws.Range("B80").Value = 0.041666666666666664
$ws.Range("C80").Value = "Koordination und Projektmanagement"
$ws.Range("D80").Value = "Teilnahme an Meeting"

# Restore the view to the top of the data (instead of scrolled down to row 20)
# and leave the selection on the newly touched helper cell.
$ws.Range("E6").Select()
